{"js": "// Add a \"Dependencies\" section at the end of the document:\n//   - a new paragraph \"Dependencies\"\n//   - a new paragraph \"Npm install -D tailwindcss postcss autoprefixer\"\n//   - a new paragraph \"Npx tailwindcss init -p\"\n// (the trailing empty paragraph that used to be last in the body is\n// replaced by these new paragraphs).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The last paragraph in the body (right before sectPr) is the empty\n// trailing paragraph the new content should take the place of.\nconst trailingParagraph = paragraphs.items[paragraphs.items.length - 1];\n\nconst dependenciesParagraph = trailingParagraph.insertParagraph(\n  \"Dependencies\",\n  \"After\"\n);\nconst npmParagraph = dependenciesParagraph.insertParagraph(\n  \"Npm install -D tailwindcss postcss autoprefixer\",\n  \"After\"\n);\nnpmParagraph.insertParagraph(\"Npx tailwindcss init -p\", \"After\");\n\n// Remove the now-redundant empty paragraph that used to be last.\ntrailingParagraph.delete();\n\nawait context.sync();\n", "ps1": "# Add a \"Dependencies\" section at the end of the document:\n#   - a new paragraph \"Dependencies\"\n#   - a new paragraph \"Npm install -D tailwindcss postcss autoprefixer\"\n#   - a new paragraph \"Npx tailwindcss init -p\"\n# (the trailing empty paragraph that used to be last in the body is\n# replaced by these new paragraphs.)\n\n$d = $word.ActiveDocument\n\n# The last paragraph in the body (right before the section break) is the\n# empty trailing paragraph the new content should take the place of.\n$trailing = $d.Paragraphs.Last\n$trailing.Range.InsertParagraphAfter()\n\n# Re-fetch after the structural edit, then type the first new paragraph.\n$d = $word.ActiveDocument\n$p1 = $d.Paragraphs.Last\n$p1.Range.InsertAfter(\"Dependencies\")\n\n$p1.Range.InsertParagraphAfter()\n$d = $word.ActiveDocument\n$p2 = $d.Paragraphs.Last\n$p2.Range.InsertAfter(\"Npm install -D tailwindcss postcss autoprefixer\")\n\n$p2.Range.InsertParagraphAfter()\n$d = $word.ActiveDocument\n$p3 = $d.Paragraphs.Last\n$p3.Range.InsertAfter(\"Npx tailwindcss init -p\")\n\n# Remove the original (now redundant) empty trailing paragraph.\n$d = $word.ActiveDocument\n$old = $d.Paragraphs.Item($d.Paragraphs.Count - 3)\n$old.Range.Delete()\n"}
